# Sync attendance_reports: swap the order of "Recorded By" names in column G
# from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}
